$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-12 Friday", 2)

$d.Content.Find.Execute("22×61=", $true, $false, $false, $false, $false, $true, 1, $false, "11×63=", 2)
$d.Content.Find.Execute("71×75=", $true, $false, $false, $false, $false, $true, 1, $false, "84×49=", 2)
$d.Content.Find.Execute("76×99=", $true, $false, $false, $false, $false, $true, 1, $false, "11×11=", 2)
$d.Content.Find.Execute("74×74=", $true, $false, $false, $false, $false, $true, 1, $false, "32×99=", 2)
$d.Content.Find.Execute("68×77=", $true, $false, $false, $false, $false, $true, 1, $false, "21×42=", 2)

$d.Content.Find.Execute("56×18=", $true, $false, $false, $false, $false, $true, 1, $false, "82×38=", 2)
$d.Content.Find.Execute("84×26=", $true, $false, $false, $false, $false, $true, 1, $false, "52×21=", 2)
$d.Content.Find.Execute("65×45=", $true, $false, $false, $false, $false, $true, 1, $false, "44×80=", 2)
$d.Content.Find.Execute("78×17=", $true, $false, $false, $false, $false, $true, 1, $false, "51×64=", 2)
$d.Content.Find.Execute("98×81=", $true, $false, $false, $false, $false, $true, 1, $false, "70×35=", 2)

$d.Content.Find.Execute("77×92=", $true, $false, $false, $false, $false, $true, 1, $false, "21×94=", 2)
$d.Content.Find.Execute("14×21=", $true, $false, $false, $false, $false, $true, 1, $false, "91×14=", 2)
$d.Content.Find.Execute("82×57=", $true, $false, $false, $false, $false, $true, 1, $false, "58×97=", 2)
$d.Content.Find.Execute("82×89=", $true, $false, $false, $false, $false, $true, 1, $false, "78×59=", 2)
$d.Content.Find.Execute("32×41=", $true, $false, $false, $false, $false, $true, 1, $false, "25×32=", 2)

$d.Content.Find.Execute("96×13=", $true, $false, $false, $false, $false, $true, 1, $false, "18×26=", 2)
$d.Content.Find.Execute("61×59=", $true, $false, $false, $false, $false, $true, 1, $false, "63×41=", 2)
$d.Content.Find.Execute("83×81=", $true, $false, $false, $false, $false, $true, 1, $false, "39×61=", 2)
$d.Content.Find.Execute("26×67=", $true, $false, $false, $false, $false, $true, 1, $false, "59×73=", 2)
$d.Content.Find.Execute("87×26=", $true, $false, $false, $false, $false, $true, 1, $false, "85×46=", 2)

$d.Content.Find.Execute("85×26=", $true, $false, $false, $false, $false, $true, 1, $false, "36×64=", 2)
$d.Content.Find.Execute("31×51=", $true, $false, $false, $false, $false, $true, 1, $false, "63×13=", 2)
$d.Content.Find.Execute("36×91=", $true, $false, $false, $false, $false, $true, 1, $false, "99×12=", 2)
$d.Content.Find.Execute("81×70=", $true, $false, $false, $false, $false, $true, 1, $false, "65×51=", 2)
$d.Content.Find.Execute("45×71=", $true, $false, $false, $false, $false, $true, 1, $false, "21×40=", 2)
